# Update column F (dSF) values for rows whose underlying data was repulled.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -2
$ws.Range("F4").Value = -10
$ws.Range("F5").Value = -3
$ws.Range("F8").Value = -6
$ws.Range("F9").Value = -5
$ws.Range("F15").Value = -2
$ws.Range("F18").Value = -9
$ws.Range("F22").Value = -4
$ws.Range("F25").Value = -3
$ws.Range("F26").Value = -14
$ws.Range("F27").Value = -3
